# Replace the "[your wish]" typo with "[optional]" on the Dockerfile
# (VS CODE) slide of the Medusa Backend Deployment deck.

$p = $ppt.ActivePresentation

$oldText = "[your wish]"
$newText = "[optional]"

$found = $false

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -like ("*" + $oldText + "*")) {
                $tr.Replace($oldText, $newText) | Out-Null
                $found = $true
            }
        }
    }
}

if (-not $found) {
    throw "Could not find text '[your wish]' to replace"
}
